# Inclusao da distancia entre apoios
# Applies the restructuring of the "Cargas pontuais" input block: replaces
# the 5-row numbered point-load table (with X/Y columns) with a simple
# label list (Carga 1..5), and inserts a new "Distancia entre os pontos de
# apoio" block above the existing "Calcular"/"Sair" controls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlEdgeLeft = 7
$xlEdgeRight = 10
$xlContinuous = 1
$xlThin = 2
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. Make room: insert two new rows at row 11 (old rows 11+ shift to 13+)
# ---------------------------------------------------------------------
$ws.Range("11:12").Insert()

# ---------------------------------------------------------------------
# 2. Row 5 header (Q(KN) / X / Y) loses its cell borders on F5/G5
# ---------------------------------------------------------------------
$ws.Range("F5").ClearFormats()
$ws.Range("G5").ClearFormats()
$ws.Range("F5").Value = "X"
$ws.Range("G5").Value = "Y"

# ---------------------------------------------------------------------
# 3. Rows 6-10: replace numbered point-load rows with label rows
# ---------------------------------------------------------------------
$ws.Range("E6:G6").ClearContents()
$ws.Range("E6:G6").ClearFormats()
$ws.Range("D6").Value = "Carga 1:"
$ws.Range("H6").ClearFormats()
$ws.Range("H6").Value = "Q (KN/m):"

$ws.Range("E7:G7").ClearContents()
$ws.Range("E7:G7").ClearFormats()
$ws.Range("D7").Value = "Carga 2:"
$ws.Range("H7").ClearFormats()
$ws.Range("H7").Value = "Início:"

$ws.Range("E8:G8").ClearContents()
$ws.Range("E8:G8").ClearFormats()
$ws.Range("D8").Value = "Carga 3:"
$ws.Range("H8").ClearFormats()
$ws.Range("H8").Value = "Fim:"

$ws.Range("E9:G9").ClearContents()
$ws.Range("E9:G9").ClearFormats()
$ws.Range("D9").Value = "Carga 4:"
$ws.Range("H9").ClearContents()
$ws.Range("H9").ClearFormats()

$ws.Range("E10:G10").ClearContents()
$ws.Range("E10:G10").ClearFormats()
$ws.Range("D10").Value = "Carga 5:"
$ws.Range("H10").ClearContents()
$ws.Range("H10").ClearFormats()
$ws.Range("I10").ClearContents()

# ---------------------------------------------------------------------
# 4. Row 11 (new): "Distancia entre os pontos de apoio em metros" banner
# ---------------------------------------------------------------------
$ws.Range("D11:H11").ClearFormats()
$ws.Range("D11:H11").ClearContents()

$ws.Range("D11").Value = "Distancia entre os pontos de apoio em metros"
$ws.Range("D11").Font.Size = 10
$ws.Range("D11").HorizontalAlignment = $xlCenter
$ws.Range("D11").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("D11").Borders.Item($xlEdgeLeft).Weight = $xlThin

$ws.Range("E11:G11").Font.Size = 10
$ws.Range("E11:G11").HorizontalAlignment = $xlCenter

$ws.Range("I11").ClearContents()
$ws.Range("I11").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I11").Borders.Item($xlEdgeRight).Weight = $xlThin

$ws.Range("D11:G11").Merge()

# ---------------------------------------------------------------------
# 5. Row 12 (new): "Distancia (m):" input + Calcular / Sair buttons
# ---------------------------------------------------------------------
$ws.Range("D12:H12").ClearFormats()
$ws.Range("D12:H12").ClearContents()

$ws.Range("D12").Value = "Distancia (m):"
$ws.Range("D12").Font.Size = 8
$ws.Range("D12").HorizontalAlignment = $xlCenter
$ws.Range("D12").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("D12").Borders.Item($xlEdgeLeft).Weight = $xlThin

$ws.Range("E12:G12").Font.Size = 10
$ws.Range("E12:G12").HorizontalAlignment = $xlCenter

$ws.Range("H12").Value = "Calcular"

$ws.Range("I12").ClearContents()
$ws.Range("I12").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I12").Borders.Item($xlEdgeRight).Weight = $xlThin
$ws.Range("I12").Value = "Sair"

# ---------------------------------------------------------------------
# 6. Selection / active cell, to match the saved workbook state
# ---------------------------------------------------------------------
$ws.Range("H9").Select()
